$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 63.91118233333333
$ws.Range("H2").Value = 191.733547
$ws.Range("I2").Value = 0.4067926910433548
$ws.Range("J2").Value = 0.4067926910433549
$ws.Range("M2").Value = 0.04532466666666666
$ws.Range("O2").Value = 0.2097030900106722
$ws.Range("P2").Value = 0.2097030900106722
$ws.Range("Q2").Value = 2.896753035530888
$ws.Range("R2").Value = 26.070777319778
$ws.Range("S2").Value = 0.0853056843055482
$ws.Range("T2").Value = 0.0853056843055482

# Row 3
$ws.Range("G3").Value = 63.91118233333333
$ws.Range("H3").Value = 191.733547
$ws.Range("I3").Value = 0.4067926910433548
$ws.Range("J3").Value = 0.4067926910433549
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.1708126666666666
$ws.Range("N3").Value = 0.512438
$ws.Range("O3").Value = 0.7902969099893278
$ws.Range("P3").Value = 0.7902969099893277
$ws.Range("Q3").Value = 10.91683948417622
$ws.Range("R3").Value = 98.25155535758599
$ws.Range("S3").Value = 0.3214870067378066
$ws.Range("T3").Value = 0.3214870067378066

# Row 4
$ws.Range("I4").Value = 0.3656254573230189
$ws.Range("J4").Value = 0.365625457323019
$ws.Range("M4").Value = 0.04532466666666666
$ws.Range("O4").Value = 0.2097030900106722
$ws.Range("P4").Value = 0.2097030900106722
$ws.Range("Q4").Value = 2.603602957199999
$ws.Range("S4").Value = 0.07667278818720223
$ws.Range("T4").Value = 0.07667278818720223

# Row 5
$ws.Range("I5").Value = 0.3656254573230189
$ws.Range("J5").Value = 0.365625457323019
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1708126666666666
$ws.Range("N5").Value = 0.512438
$ws.Range("O5").Value = 0.7902969099893278
$ws.Range("P5").Value = 0.7902969099893277
$ws.Range("Q5").Value = 9.812060336399998
$ws.Range("R5").Value = 88.30854302759998
$ws.Range("S5").Value = 0.2889526691358167
$ws.Range("T5").Value = 0.2889526691358167

# Row 6
$ws.Range("G6").Value = 35.755375
$ws.Range("H6").Value = 107.266125
$ws.Range("I6").Value = 0.2275818516336261
$ws.Range("J6").Value = 0.2275818516336262
$ws.Range("M6").Value = 0.04532466666666666
$ws.Range("O6").Value = 0.2097030900106722
$ws.Range("P6").Value = 0.2097030900106722
$ws.Range("Q6").Value = 1.620600453416666
$ws.Range("R6").Value = 14.58540408075
$ws.Range("S6").Value = 0.04772461751792175
$ws.Range("T6").Value = 0.04772461751792175

# Row 7
$ws.Range("G7").Value = 35.755375
$ws.Range("H7").Value = 107.266125
$ws.Range("I7").Value = 0.2275818516336261
$ws.Range("J7").Value = 0.2275818516336262
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.1708126666666666
$ws.Range("N7").Value = 0.512438
$ws.Range("O7").Value = 0.7902969099893278
$ws.Range("P7").Value = 0.7902969099893277
$ws.Range("Q7").Value = 6.107470951416666
$ws.Range("R7").Value = 54.96723856275
$ws.Range("S7").Value = 0.1798572341157044
$ws.Range("T7").Value = 0.1798572341157044
